$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the segment names currently stored in column A before we shift things.
$names = @()
For ($i = 0; $i -le 18; $i++) {
    $names += $ws.Cells.Item(2 + $i, 1).Value2
}

# Insert a new column before column B. This shifts the old B (PercActivations)
# and C (PercSegmentAreas) columns one to the right (-> C, D). Column A (with
# the segment names / its header style) is untouched by this operation.
$ws.Columns("B").Insert()

# Header row.
$ws.Range("B1").Value = "segments"
# New B1 header needs the same bordered/bold header formatting as the other
# header cells (copy it from C1, which already carries that style).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# Column A now holds the 0-based numeric index (keeps its original bordered
# header-ish style); column B holds the segment name text (moved out of
# column A) and should be a plain, unstyled cell like the values that used
# to live in B/C.
For ($i = 0; $i -le 18; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 1).Value = $i
}

$ws.Range("C2:C20").Copy()
$ws.Range("B2:B20").PasteSpecial(-4122)
